# Reorders the elements (lists of node ids) contained in certain cell
# values across the workbook's sheets, matching the canonical OOXML diff.
# Only the ordering of the list items inside each cell's string changes;
# the set of items stays the same.

$wb = $excel.ActiveWorkbook

$wsU18 = $wb.Worksheets.Item("Node U18")
$wsU18.Range("G3").Value = "['U109', 'U124', 'U54', 'U3', 'U79']"
$wsU18.Range("G5").Value = "['U109', 'U54', 'U3', 'U79', 'U99']"
$wsU18.Range("I6").Value = "['U130', 'U62']"

$wsU142 = $wb.Worksheets.Item("Node U142")
$wsU142.Range("E3").Value = "['U130', 'U54', 'U47', 'U67', 'U79', 'U91', 'U110']"
$wsU142.Range("G3").Value = "['U42', 'U123', 'U10']"
$wsU142.Range("G4").Value = "['U42', 'U10']"
$wsU142.Range("I6").Value = "['U68', 'U4']"

$wsU42 = $wb.Worksheets.Item("Node U42")
$wsU42.Range("E4").Value = "['U106', 'U47', 'U118', 'U41']"
$wsU42.Range("G2").Value = "['U123', 'U142']"

$wsU23 = $wb.Worksheets.Item("Node U23")
$wsU23.Range("I2").Value = "['U19', 'U73', 'U17', 'U14', 'U1']"
$wsU23.Range("I3").Value = "['U19', 'U73', 'U17', 'U14', 'U1']"
$wsU23.Range("I4").Value = "['U19', 'U73', 'U17', 'U14', 'U1']"

$wsU90 = $wb.Worksheets.Item("Node U90")
$wsU90.Range("G2").Value = "['U109', 'U79', 'U126', 'U76', 'U6']"
$wsU90.Range("I2").Value = "['U54', 'U3']"
$wsU90.Range("E3").Value = "['U134', 'U62', 'U18']"
$wsU90.Range("G3").Value = "['U109', 'U79', 'U126', 'U76', 'U6']"
$wsU90.Range("I3").Value = "['U54', 'U3']"
$wsU90.Range("E4").Value = "['U123', 'U4', 'U118']"
$wsU90.Range("I4").Value = "['U54', 'U3']"
